# Update Betfair odds values for Jogos_do_Dia_Betfair_Back_Lay_2025-12-31.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 1.37
$ws.Range("S2").Value = 3.75
$ws.Range("AH2").Value = 18.5

# Row 6
$ws.Range("G6").Value = 8.6
$ws.Range("H6").Value = 1.46
$ws.Range("N6").Value = 4.4
$ws.Range("S6").Value = 2.42
$ws.Range("T6").Value = 1.79
$ws.Range("V6").Value = 2.68
$ws.Range("W6").Value = 1.15
$ws.Range("AF6").Value = 70

# Row 7
$ws.Range("O7").Value = 1.16
$ws.Range("Q7").Value = 1.5
$ws.Range("AB7").Value = 22
$ws.Range("AD7").Value = 12.5
$ws.Range("AE7").Value = 21
$ws.Range("AJ7").Value = 60
$ws.Range("AK7").Value = 34
$ws.Range("AL7").Value = 36
$ws.Range("AN7").Value = 21

# Row 8
$ws.Range("F8").Value = 2.08
$ws.Range("G8").Value = 2.38
$ws.Range("I8").Value = 4.1
$ws.Range("J8").Value = 3.4

# Row 9
$ws.Range("Z9").Value = 1000

# Row 10
$ws.Range("W10").Value = 1.01

# Row 12
$ws.Range("G12").Value = 1.36
$ws.Range("L12").Value = 1.19
$ws.Range("W12").Value = 3.55
$ws.Range("AG12").Value = 13.5
